$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11..19 down to 12..20
$ws.Rows("11:11").Insert()

# Populate the new row 11 with the newest weekly data entry
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44467
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100114007
$ws.Range("G11").Value = "Jengibre"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("N11").Value = "`$/caja 13 kilos"
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 1038
$ws.Range("Q11").Value = 13
$ws.Range("R11").Value = "Hortaliza"
